function Get-ShapeByName {
    param($Slide, [string]$Name)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $sh = $Slide.Shapes.Item($i)
        if ($sh.Name -eq $Name) { return $sh }
    }
    return $Slide.Shapes.Item(1)
}

$p = $ppt.ActivePresentation

# --- Slide 1: Title "Scriptory" -> "Timelink" -----------------------------
$slide1 = $p.Slides.Item(1)
$titleShape = Get-ShapeByName $slide1 "Title 1"
$titleRange = $titleShape.TextFrame.TextRange
# "Scriptory" is the first 9 characters of the title text run; replacing just
# that run of characters keeps the following line breaks / "Browser game"
# run untouched.
$titleRange.Characters(1, 9).Text = "Timelink"

# --- Slide 2: Content placeholder wording tweaks --------------------------
$slide2 = $p.Slides.Item(2)
$bodyShape = Get-ShapeByName $slide2 "Content Placeholder 2"
$bodyRange = $bodyShape.TextFrame.TextRange

# Paragraph 2: drop the "(the Timeline)" aside.
$para2 = $bodyRange.Paragraphs(2, 1)
$para2.Runs(1, 1).Text = "The playground field consists of different historical ages, which user needs to pass through to win."

# Paragraph 3: rename "Timeline" -> "Timelink" (this becomes its own run,
# split out of the sentence it lives in).
$para3 = $bodyRange.Paragraphs(3, 1)
$null = $para3.Replace("Timeline", "Timelink", 0, $false, $false)
